# Applies targeted odds/value updates to Sheet1 cells per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 7).Value = 2.55  # G3
$ws.Cells.Item(3, 9).Value = 3.3  # I3
$ws.Cells.Item(3, 10).Value = 3.5  # J3
$ws.Cells.Item(3, 12).Value = 4  # L3
$ws.Cells.Item(3, 15).Value = 1.62  # O3
$ws.Cells.Item(3, 16).Value = 2.2  # P3
$ws.Cells.Item(3, 19).Value = 2.88  # S3
$ws.Cells.Item(3, 20).Value = 1.4  # T3
$ws.Cells.Item(3, 23).Value = 6.5  # W3
$ws.Cells.Item(3, 24).Value = 1.11  # X3
$ws.Cells.Item(3, 32).Value = 26  # AF3
$ws.Cells.Item(3, 41).Value = 13  # AO3
$ws.Cells.Item(3, 43).Value = 34  # AQ3

# Row 4
$ws.Cells.Item(4, 7).Value = 2.5  # G4
$ws.Cells.Item(4, 9).Value = 3.4  # I4
$ws.Cells.Item(4, 10).Value = 3.5  # J4
$ws.Cells.Item(4, 20).Value = 1.33  # T4
$ws.Cells.Item(4, 27).Value = 2.38  # AA4
$ws.Cells.Item(4, 28).Value = 1.53  # AB4
$ws.Cells.Item(4, 30).Value = 10  # AD4
$ws.Cells.Item(4, 32).Value = 26  # AF4
$ws.Cells.Item(4, 40).Value = 6.5  # AN4

# Row 5
$ws.Cells.Item(5, 7).Value = 1.73  # G5
$ws.Cells.Item(5, 12).Value = 7.5  # L5
$ws.Cells.Item(5, 29).Value = 4.33  # AC5
$ws.Cells.Item(5, 32).Value = 13  # AF5
$ws.Cells.Item(5, 36).Value = 7.5  # AJ5
$ws.Cells.Item(5, 43).Value = 67  # AQ5

# Row 9
$ws.Cells.Item(9, 7).Value = 1.42  # G9
$ws.Cells.Item(9, 8).Value = 4.1  # H9
$ws.Cells.Item(9, 9).Value = 7.8  # I9
$ws.Cells.Item(9, 10).Value = 1.95  # J9
$ws.Cells.Item(9, 11).Value = 2.22  # K9
$ws.Cells.Item(9, 12).Value = 6.9  # L9
$ws.Cells.Item(9, 14).Value = 7.6  # N9
$ws.Cells.Item(9, 15).Value = 1.27  # O9
$ws.Cells.Item(9, 16).Value = 3.45  # P9
$ws.Cells.Item(9, 19).Value = 1.8  # S9
$ws.Cells.Item(9, 20).Value = 1.91  # T9
$ws.Cells.Item(9, 23).Value = 2.9  # W9
$ws.Cells.Item(9, 24).Value = 1.36  # X9
$ws.Cells.Item(9, 25).Value = 1.39  # Y9
$ws.Cells.Item(9, 26).Value = 2.77  # Z9
$ws.Cells.Item(9, 30).Value = 6.2  # AD9
$ws.Cells.Item(9, 32).Value = 9  # AF9
$ws.Cells.Item(9, 33).Value = 12  # AG9
$ws.Cells.Item(9, 34).Value = 29  # AH9
$ws.Cells.Item(9, 35).Value = 7.6  # AI9
$ws.Cells.Item(9, 36).Value = 8.25  # AJ9
$ws.Cells.Item(9, 37).Value = 19.5  # AK9
$ws.Cells.Item(9, 40).Value = 19.5  # AN9
$ws.Cells.Item(9, 41).Value = 55  # AO9
$ws.Cells.Item(9, 42).Value = 23  # AP9
$ws.Cells.Item(9, 43).Value = 200  # AQ9
$ws.Cells.Item(9, 44).Value = 90  # AR9
$ws.Cells.Item(9, 45).Value = 75  # AS9

# Row 12
$ws.Cells.Item(12, 7).Value = 1.14  # G12
$ws.Cells.Item(12, 8).Value = 7  # H12
$ws.Cells.Item(12, 9).Value = 12  # I12
$ws.Cells.Item(12, 10).Value = 1.5  # J12
$ws.Cells.Item(12, 11).Value = 2.88  # K12
$ws.Cells.Item(12, 12).Value = 10  # L12
$ws.Cells.Item(12, 13).Value = 1.02  # M12
$ws.Cells.Item(12, 14).Value = 12  # N12
$ws.Cells.Item(12, 15).Value = 1.13  # O12
$ws.Cells.Item(12, 16).Value = 5.5  # P12
$ws.Cells.Item(12, 25).Value = 1.22  # Y12
$ws.Cells.Item(12, 26).Value = 4  # Z12
$ws.Cells.Item(12, 27).Value = 2.2  # AA12
$ws.Cells.Item(12, 28).Value = 1.62  # AB12
$ws.Cells.Item(12, 29).Value = 9.5  # AC12
$ws.Cells.Item(12, 31).Value = 11  # AE12
$ws.Cells.Item(12, 32).Value = 7  # AF12
$ws.Cells.Item(12, 34).Value = 34  # AH12
$ws.Cells.Item(12, 35).Value = 19  # AI12
$ws.Cells.Item(12, 36).Value = 15  # AJ12
$ws.Cells.Item(12, 37).Value = 29  # AK12
$ws.Cells.Item(12, 38).Value = 81  # AL12
$ws.Cells.Item(12, 39).Value = 1250  # AM12
$ws.Cells.Item(12, 40).Value = 34  # AN12
$ws.Cells.Item(12, 42).Value = 34  # AP12
$ws.Cells.Item(12, 43).Value = 151  # AQ12
$ws.Cells.Item(12, 44).Value = 81  # AR12
$ws.Cells.Item(12, 45).Value = 67  # AS12

# Row 13
$ws.Cells.Item(13, 13).Value = 1.01  # M13
$ws.Cells.Item(13, 14).Value = 17  # N13
$ws.Cells.Item(13, 21).Value = 1.93  # U13
$ws.Cells.Item(13, 22).Value = 1.88  # V13
$ws.Cells.Item(13, 27).Value = 1.62  # AA13
$ws.Cells.Item(13, 32).Value = 15  # AF13

# Row 14
$ws.Cells.Item(14, 7).Value = 5.5  # G14
$ws.Cells.Item(14, 8).Value = 4.33  # H14
$ws.Cells.Item(14, 9).Value = 1.48  # I14
$ws.Cells.Item(14, 10).Value = 6  # J14
$ws.Cells.Item(14, 11).Value = 2.3  # K14
$ws.Cells.Item(14, 12).Value = 2  # L14
$ws.Cells.Item(14, 15).Value = 1.22  # O14
$ws.Cells.Item(14, 16).Value = 4  # P14
$ws.Cells.Item(14, 19).Value = 1.73  # S14
$ws.Cells.Item(14, 20).Value = 2.08  # T14
$ws.Cells.Item(14, 23).Value = 2.75  # W14
$ws.Cells.Item(14, 24).Value = 1.4  # X14
$ws.Cells.Item(14, 27).Value = 1.91  # AA14
$ws.Cells.Item(14, 29).Value = 15  # AC14
$ws.Cells.Item(14, 30).Value = 29  # AD14
$ws.Cells.Item(14, 31).Value = 19  # AE14
$ws.Cells.Item(14, 32).Value = 67  # AF14
$ws.Cells.Item(14, 36).Value = 8.5  # AJ14
$ws.Cells.Item(14, 37).Value = 19  # AK14
$ws.Cells.Item(14, 39).Value = 600  # AM14
$ws.Cells.Item(14, 41).Value = 7  # AO14
$ws.Cells.Item(14, 43).Value = 10  # AQ14

# Row 15
$ws.Cells.Item(15, 13).Value = 1.03  # M15
$ws.Cells.Item(15, 14).Value = 17  # N15
$ws.Cells.Item(15, 15).Value = 1.14  # O15
$ws.Cells.Item(15, 21).Value = 1.8  # U15
$ws.Cells.Item(15, 22).Value = 2.05  # V15
$ws.Cells.Item(15, 24).Value = 1.62  # X15
$ws.Cells.Item(15, 25).Value = 1.29  # Y15
$ws.Cells.Item(15, 26).Value = 3.5  # Z15
$ws.Cells.Item(15, 27).Value = 1.53  # AA15
$ws.Cells.Item(15, 28).Value = 2.38  # AB15
$ws.Cells.Item(15, 29).Value = 10  # AC15
$ws.Cells.Item(15, 30).Value = 10  # AD15
$ws.Cells.Item(15, 35).Value = 17  # AI15

# Row 16
$ws.Cells.Item(16, 15).Value = 1.22  # O16
$ws.Cells.Item(16, 16).Value = 4  # P16
$ws.Cells.Item(16, 19).Value = 1.8  # S16
$ws.Cells.Item(16, 20).Value = 2  # T16
$ws.Cells.Item(16, 23).Value = 2.75  # W16
$ws.Cells.Item(16, 24).Value = 1.4  # X16
$ws.Cells.Item(16, 27).Value = 1.67  # AA16

# Row 17
$ws.Cells.Item(17, 8).Value = 4.1  # H17
$ws.Cells.Item(17, 10).Value = 5.1  # J17
$ws.Cells.Item(17, 12).Value = 2.05  # L17
$ws.Cells.Item(17, 16).Value = 3.85  # P17
$ws.Cells.Item(17, 25).Value = 1.32  # Y17
$ws.Cells.Item(17, 26).Value = 3.1  # Z17
$ws.Cells.Item(17, 29).Value = 15.5  # AC17
$ws.Cells.Item(17, 33).Value = 50  # AG17
$ws.Cells.Item(17, 34).Value = 50  # AH17
$ws.Cells.Item(17, 41).Value = 7.9  # AO17

# Row 18
$ws.Cells.Item(18, 7).Value = 1.38  # G18
$ws.Cells.Item(18, 8).Value = 4.4  # H18
$ws.Cells.Item(18, 9).Value = 7.2  # I18
$ws.Cells.Item(18, 10).Value = 1.85  # J18
$ws.Cells.Item(18, 11).Value = 2.32  # K18
$ws.Cells.Item(18, 12).Value = 6.6  # L18
$ws.Cells.Item(18, 15).Value = 1.25  # O18
$ws.Cells.Item(18, 16).Value = 3.25  # P18
$ws.Cells.Item(18, 19).Value = 1.75  # S18
$ws.Cells.Item(18, 20).Value = 1.87  # T18
$ws.Cells.Item(18, 23).Value = 2.72  # W18
$ws.Cells.Item(18, 24).Value = 1.34  # X18
$ws.Cells.Item(18, 29).Value = 6.3  # AC18
$ws.Cells.Item(18, 32).Value = 8.5  # AF18
$ws.Cells.Item(18, 35).Value = 11  # AI18
$ws.Cells.Item(18, 36).Value = 8.75  # AJ18
$ws.Cells.Item(18, 37).Value = 23  # AK18
$ws.Cells.Item(18, 40).Value = 17.5  # AN18
$ws.Cells.Item(18, 42).Value = 24  # AP18
$ws.Cells.Item(18, 43).Value = 175  # AQ18
$ws.Cells.Item(18, 44).Value = 90  # AR18
$ws.Cells.Item(18, 45).Value = 90  # AS18
